$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values down to plain integers
$ws.Range("Q2").Value = 510540
$ws.Range("R2").Value = 6544167

# Remove the Starttid (Z2) and Sluttid (AB2) values entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
